# Timesheet_Group5.xlsx - update by sravani date:23/02/2013
#
# For rows 32-39 on the "February 2013" sheet, the "OFF" marker that used to
# sit in column AA is moved one column to the right (column AC), and the two
# now-vacated cells (AA, AB) receive numeric hour values. Columns AD..AH are
# left untouched.
#
# Row -> AA value, AB value
#  32 -> 0   , 0
#  33 -> 0   , 0
#  34 -> 0   , 0
#  35 -> 1.5 , 1
#  36 -> 0   , 0
#  37 -> 0   , 0
#  38 -> 0   , 0
#  39 -> 1.5 , 1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("February 2013")

$rowValues = @{
    32 = @(0, 0)
    33 = @(0, 0)
    34 = @(0, 0)
    35 = @(1.5, 1)
    36 = @(0, 0)
    37 = @(0, 0)
    38 = @(0, 0)
    39 = @(1.5, 1)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]

    # AA/AB currently hold "OFF" (AA) and a blank (AB); restyle them to the
    # plain numeric-cell look (same format already used by the neighboring
    # blank cell AD) before writing numbers into them.
    $ws.Range("AD$r").Copy()
    $ws.Range("AA$r" + ":" + "AB$r").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Range("AA$r").Value = $vals[0]
    $ws.Range("AB$r").Value = $vals[1]

    # AC becomes the new "OFF" marker cell - copy the grey "OFF" formatting
    # from Z (still "OFF" on every row) then set its text.
    $ws.Range("Z$r").Copy()
    $ws.Range("AC$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("AC$r").Value = "OFF"
}

$excel.CutCopyMode = $false

# Move the active selection to AA35, matching the saved view state.
$ws.Range("AA35").Select()
